$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume Number + report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/24/2024  Through  6/30/2024"

# --- Cells that change data TYPE/STYLE (number <-> text) ---
# Copy style (and base content) from a stable same-style neighbor first,
# then overwrite with the correct final value/text so the OOXML cell
# ends up with the right s= (style) and t= (type) attributes.
$ws.Range("F22").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 1

$ws.Range("C14").Copy($ws.Range("C23"))

$ws.Range("C14").Copy($ws.Range("F29"))
$ws.Range("C14").Copy($ws.Range("G29"))
$ws.Range("E14").Copy($ws.Range("H29"))

$ws.Range("C14").Copy($ws.Range("F30"))
$ws.Range("C14").Copy($ws.Range("G30"))
$ws.Range("E14").Copy($ws.Range("H30"))

# --- Plain numeric updates (weekly crime-stat refresh) ---
$ws.Range("M14").Value = -83.333333333333
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = -66.666666666666
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -60
$ws.Range("I15").Value = 16
$ws.Range("J15").Value = 13
$ws.Range("K15").Value = 23.076923076923
$ws.Range("L15").Value = -5.882352941176
$ws.Range("N15").Value = -55.555555555555
$ws.Range("C16").Value = 11
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 266.666666666667
$ws.Range("F16").Value = 33
$ws.Range("G16").Value = 27
$ws.Range("H16").Value = 22.222222222222
$ws.Range("I16").Value = 188
$ws.Range("J16").Value = 158
$ws.Range("K16").Value = 18.987341772151
$ws.Range("L16").Value = 33.333333333333
$ws.Range("M16").Value = -3.092783505154
$ws.Range("N16").Value = -73.925104022191
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 50
$ws.Range("G17").Value = 36
$ws.Range("H17").Value = -16.666666666666
$ws.Range("I17").Value = 213
$ws.Range("J17").Value = 218
$ws.Range("K17").Value = -2.293577981651
$ws.Range("L17").Value = 29.090909090909
$ws.Range("M17").Value = 18.994413407821
$ws.Range("N17").Value = -49.043062200956
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -75
$ws.Range("I18").Value = 110
$ws.Range("J18").Value = 110
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = -25.675675675675
$ws.Range("M18").Value = -34.523809523809
$ws.Range("N18").Value = -82.456140350877
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 22
$ws.Range("E19").Value = -63.636363636363
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 65
$ws.Range("H19").Value = -41.538461538461
$ws.Range("I19").Value = 292
$ws.Range("J19").Value = 355
$ws.Range("K19").Value = -17.746478873239
$ws.Range("L19").Value = -9.316770186335
$ws.Range("M19").Value = 105.633802816901
$ws.Range("N19").Value = 8.550185873605
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 150
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 25
$ws.Range("H20").Value = -40
$ws.Range("I20").Value = 100
$ws.Range("J20").Value = 96
$ws.Range("K20").Value = 4.166666666666
$ws.Range("L20").Value = -5.660377358490
$ws.Range("M20").Value = 28.205128205128
$ws.Range("N20").Value = -79.466119096509
$ws.Range("C21").Value = 38
$ws.Range("D21").Value = 41
$ws.Range("E21").Value = -7.317073170731
$ws.Range("F21").Value = 124
$ws.Range("G21").Value = 183
$ws.Range("H21").Value = -32.240437158469
$ws.Range("I21").Value = 920
$ws.Range("J21").Value = 952
$ws.Range("K21").Value = -3.361344537815
$ws.Range("L21").Value = 2.108768035516
$ws.Range("M21").Value = 17.496807151979
$ws.Range("N21").Value = -64.271844660194
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -50
$ws.Range("I22").Value = 18
$ws.Range("J22").Value = 17
$ws.Range("K22").Value = 5.882352941176
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = 80
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -100
$ws.Range("J23").Value = 22
$ws.Range("K23").Value = -40.909090909090
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -14.285714285714
$ws.Range("F24").Value = 80
$ws.Range("G24").Value = 112
$ws.Range("H24").Value = -28.571428571428
$ws.Range("I24").Value = 468
$ws.Range("J24").Value = 512
$ws.Range("K24").Value = -8.59375
$ws.Range("L24").Value = -2.904564315352
$ws.Range("M24").Value = 31.830985915493
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 133.333333333333
$ws.Range("F25").Value = 24
$ws.Range("G25").Value = 11
$ws.Range("H25").Value = 118.181818181818
$ws.Range("I25").Value = 118
$ws.Range("J25").Value = 81
$ws.Range("K25").Value = 45.679012345679
$ws.Range("L25").Value = 34.090909090909
$ws.Range("C26").Value = 21
$ws.Range("D26").Value = 12
$ws.Range("E26").Value = 75
$ws.Range("F26").Value = 63
$ws.Range("G26").Value = 52
$ws.Range("H26").Value = 21.153846153846
$ws.Range("I26").Value = 338
$ws.Range("J26").Value = 337
$ws.Range("K26").Value = 0.296735905044
$ws.Range("L26").Value = 9.032258064516
$ws.Range("M26").Value = -13.994910941475
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -66.666666666666
$ws.Range("I27").Value = 24
$ws.Range("J27").Value = 23
$ws.Range("K27").Value = 4.347826086956
$ws.Range("L27").Value = 9.090909090909
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 7
$ws.Range("H28").Value = 0
$ws.Range("L28").Value = 30
$ws.Range("M29").Value = -90.909090909090
$ws.Range("N29").Value = -98.039215686274
$ws.Range("M30").Value = -87.5
$ws.Range("N30").Value = -97.959183673469
